# NATMI ligand-receptor export (Lama3-Sdc2): refresh the expression /
# specificity / edge-weight figures for the "ECs" cluster using the
# updated TPM values (commit: "update scripts wuth new tpm").
#
# Only the "ECs" sending/target-cluster TPM figures moved; the per-cell
# ripple below is the recomputed Ligand/Receptor average+total expression,
# derived specificity (share of the 3-cluster total), and the edge
# weights/specificities (products of the ligand & receptor figures) for
# every row that references "ECs" as sender and/or target, plus the
# renormalised specificity of the two rows that don't.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.505203
$ws.Range("H2").Value = 10.515609
$ws.Range("I2").Value = 0.445953142814298
$ws.Range("J2").Value = 0.445953142814298
$ws.Range("M2").Value = 1.116695
$ws.Range("N2").Value = 3.350085
$ws.Range("O2").Value = 0.008174214292497491
$ws.Range("P2").Value = 0.008174214292497492
$ws.Range("Q2").Value = 3.914242664085
$ws.Range("R2").Value = 35.228183976765
$ws.Range("S2").Value = 0.003645316553776809
$ws.Range("T2").Value = 0.00364531655377681
$ws.Range("G3").Value = 3.505203
$ws.Range("H3").Value = 10.515609
$ws.Range("I3").Value = 0.445953142814298
$ws.Range("J3").Value = 0.445953142814298
$ws.Range("O3").Value = 0.8193429796700005
$ws.Range("P3").Value = 0.8193429796700005
$ws.Range("Q3").Value = 392.34440556402
$ws.Range("R3").Value = 3531.09965007618
$ws.Range("S3").Value = 0.3653885768266682
$ws.Range("T3").Value = 0.3653885768266682
$ws.Range("G4").Value = 3.505203
$ws.Range("H4").Value = 10.515609
$ws.Range("I4").Value = 0.445953142814298
$ws.Range("J4").Value = 0.445953142814298
$ws.Range("O4").Value = 0.172482806037502
$ws.Range("P4").Value = 0.1724828060375021
$ws.Range("Q4").Value = 82.593816855615
$ws.Range("R4").Value = 743.344351700535
$ws.Range("S4").Value = 0.076919249433853
$ws.Range("T4").Value = 0.07691924943385302
$ws.Range("I5").Value = 0.01738179300185462
$ws.Range("J5").Value = 0.01738179300185462
$ws.Range("M5").Value = 1.116695
$ws.Range("N5").Value = 3.350085
$ws.Range("O5").Value = 0.008174214292497491
$ws.Range("P5").Value = 0.008174214292497492
$ws.Range("Q5").Value = 0.1525643598266667
$ws.Range("R5").Value = 1.37307923844
$ws.Range("S5").Value = 0.0001420825007849929
$ws.Range("T5").Value = 0.0001420825007849929
$ws.Range("I6").Value = 0.01738179300185462
$ws.Range("J6").Value = 0.01738179300185462
$ws.Range("O6").Value = 0.8193429796700005
$ws.Range("P6").Value = 0.8193429796700005
$ws.Range("S6").Value = 0.01424165007014673
$ws.Range("T6").Value = 0.01424165007014672
$ws.Range("I7").Value = 0.01738179300185462
$ws.Range("J7").Value = 0.01738179300185462
$ws.Range("O7").Value = 0.172482806037502
$ws.Range("P7").Value = 0.1724828060375021
$ws.Range("S7").Value = 0.002998060430922901
$ws.Range("T7").Value = 0.002998060430922901
$ws.Range("I8").Value = 0.5366650641838474
$ws.Range("J8").Value = 0.5366650641838474
$ws.Range("M8").Value = 1.116695
$ws.Range("N8").Value = 3.350085
$ws.Range("O8").Value = 0.008174214292497491
$ws.Range("P8").Value = 0.008174214292497492
$ws.Range("Q8").Value = 4.710443965694999
$ws.Range("R8").Value = 42.393995691255
$ws.Range("S8").Value = 0.004386815237935689
$ws.Range("T8").Value = 0.00438681523793569
$ws.Range("I9").Value = 0.5366650641838474
$ws.Range("J9").Value = 0.5366650641838474
$ws.Range("O9").Value = 0.8193429796700005
$ws.Range("P9").Value = 0.8193429796700005
$ws.Range("S9").Value = 0.4397127527731856
$ws.Range("T9").Value = 0.4397127527731856
$ws.Range("I10").Value = 0.5366650641838474
$ws.Range("J10").Value = 0.5366650641838474
$ws.Range("O10").Value = 0.172482806037502
$ws.Range("P10").Value = 0.1724828060375021
$ws.Range("S10").Value = 0.09256549617272614
$ws.Range("T10").Value = 0.09256549617272615
